# Auto-generated edit script applying numeric updates to Behemoth_Profits workbook
# (columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#  K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 715.3570999999999
$ws.Range("J17").Value = 715.3570999999999
$ws.Range("L17").Value = 2146.0713
$ws.Range("N17").Value = -2482.0713
$ws.Range("H100").Value = 1985.2858
$ws.Range("J100").Value = 2185.4
$ws.Range("L100").Value = 2185.4
$ws.Range("N100").Value = -3267.4
$ws.Range("H106").Value = 7519.353
$ws.Range("I106").Value = 1738.9375
$ws.Range("K106").Value = 1738.9375
$ws.Range("M106").Value = -1107.9375
$ws.Range("H112").Value = 2118.625
$ws.Range("J112").Value = 2118.625
$ws.Range("L112").Value = 6355.875
$ws.Range("N112").Value = -8571.875
$ws.Range("H132").Value = 1535.7709
$ws.Range("I132").Value = 1392.7222
$ws.Range("J132").Value = 1964.9166
$ws.Range("K132").Value = 4178.1666
$ws.Range("L132").Value = 5894.7498
$ws.Range("M132").Value = -1648.1666
$ws.Range("N132").Value = -10954.7498
$ws.Range("H137").Value = 14672.777
$ws.Range("I137").Value = 32542.348
$ws.Range("K137").Value = 97627.04400000001
$ws.Range("M137").Value = -95077.04400000001
$ws.Range("H138").Value = 4607.9033
$ws.Range("J138").Value = 5473.05
$ws.Range("L138").Value = 16419.15
$ws.Range("N138").Value = -26699.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7585283
$ws.Range("I32").Value = 7943264
$ws.Range("K32").Value = 7943264
$ws.Range("M32").Value = -7942977
$ws.Range("H50").Value = 6719
$ws.Range("I50").Value = 6387.4
$ws.Range("J50").Value = 7548
$ws.Range("K50").Value = 6387.4
$ws.Range("L50").Value = 7548
$ws.Range("M50").Value = -5673.4
$ws.Range("N50").Value = -8976
$ws.Range("H54").Value = 39495
$ws.Range("J54").Value = 39495
$ws.Range("L54").Value = 39495
$ws.Range("N54").Value = -41033
$ws.Range("H60").Value = 15333.333
$ws.Range("J60").Value = 6000
$ws.Range("L60").Value = 6000
$ws.Range("N60").Value = -7466
$ws.Range("H74").Value = 10008161
$ws.Range("I74").Value = 14707769
$ws.Range("J74").Value = 21495.125
$ws.Range("K74").Value = 14707769
$ws.Range("L74").Value = 21495.125
$ws.Range("M74").Value = -14706895
$ws.Range("N74").Value = -23243.125
$ws.Range("H77").Value = 10008161
$ws.Range("I77").Value = 14707769
$ws.Range("J77").Value = 21495.125
$ws.Range("K77").Value = 73538845
$ws.Range("L77").Value = 107475.625
$ws.Range("M77").Value = -73534477
$ws.Range("N77").Value = -116211.625
$ws.Range("H101").Value = 85708
$ws.Range("J101").Value = 85708
$ws.Range("L101").Value = 85708
$ws.Range("N101").Value = -92198
$ws.Range("H122").Value = 3693.682
$ws.Range("I122").Value = 1686.375
$ws.Range("K122").Value = 5059.125
$ws.Range("M122").Value = -2609.125
$ws.Range("H132").Value = 4938.95
$ws.Range("I132").Value = 1662.6154
$ws.Range("J132").Value = 11023.571
$ws.Range("K132").Value = 4987.8462
$ws.Range("L132").Value = 33070.713
$ws.Range("M132").Value = -2457.8462
$ws.Range("N132").Value = -38130.713

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1871.8889
$ws.Range("I86").Value = 1507.8334
$ws.Range("K86").Value = 1507.8334
$ws.Range("M86").Value = -384.8334
$ws.Range("H89").Value = 1871.8889
$ws.Range("I89").Value = 1507.8334
$ws.Range("K89").Value = 7539.166999999999
$ws.Range("M89").Value = -1923.166999999999
$ws.Range("H122").Value = 98926.664
$ws.Range("J122").Value = 98926.664
$ws.Range("L122").Value = 98926.664
$ws.Range("N122").Value = -108726.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 523403.7
$ws.Range("I31").Value = 8464.348
$ws.Range("K31").Value = 8464.348
$ws.Range("M31").Value = -8169.348
$ws.Range("H34").Value = 523403.7
$ws.Range("I34").Value = 8464.348
$ws.Range("K34").Value = 8464.348
$ws.Range("M34").Value = -8262.348
$ws.Range("H58").Value = 2705
$ws.Range("I58").Value = 1661.6
$ws.Range("J58").Value = 4444
$ws.Range("K58").Value = 1661.6
$ws.Range("L58").Value = 4444
$ws.Range("M58").Value = -1458.6
$ws.Range("N58").Value = -4850
$ws.Range("H122").Value = 4095.4092
$ws.Range("I122").Value = 2873.4614
$ws.Range("K122").Value = 8620.3842
$ws.Range("M122").Value = -6170.3842
$ws.Range("H132").Value = 2951.0527
$ws.Range("I132").Value = 2711
$ws.Range("K132").Value = 8133
$ws.Range("M132").Value = -5603
$ws.Range("H136").Value = 2705
$ws.Range("I136").Value = 1661.6
$ws.Range("J136").Value = 4444
$ws.Range("K136").Value = 4984.799999999999
$ws.Range("L136").Value = 13332
$ws.Range("M136").Value = -2434.799999999999
$ws.Range("N136").Value = -18432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1022.375
$ws.Range("I32").Value = 999.75
$ws.Range("J32").Value = 1045
$ws.Range("K32").Value = 2999.25
$ws.Range("L32").Value = 3135
$ws.Range("M32").Value = -2716.25
$ws.Range("N32").Value = -3701
$ws.Range("H70").Value = 9000
$ws.Range("J70").Value = 9000
$ws.Range("L70").Value = 27000
$ws.Range("N70").Value = -27630
$ws.Range("H73").Value = 9000
$ws.Range("J73").Value = 9000
$ws.Range("L73").Value = 27000
$ws.Range("N73").Value = -29184
$ws.Range("H86").Value = 2498.0908
$ws.Range("I86").Value = 619.75
$ws.Range("K86").Value = 1859.25
$ws.Range("M86").Value = -673.25
$ws.Range("H89").Value = 2498.0908
$ws.Range("I89").Value = 619.75
$ws.Range("K89").Value = 5577.75
$ws.Range("M89").Value = 350.25
$ws.Range("H107").Value = 629.2
$ws.Range("I107").Value = 582.0833
$ws.Range("J107").Value = 699.875
$ws.Range("K107").Value = 1746.2499
$ws.Range("L107").Value = 2099.625
$ws.Range("M107").Value = 173.7501
$ws.Range("N107").Value = -5939.625
$ws.Range("H113").Value = 1134.3636
$ws.Range("I113").Value = 537.2143
$ws.Range("J113").Value = 1574.3684
$ws.Range("K113").Value = 1611.6429
$ws.Range("L113").Value = 4723.1052
$ws.Range("M113").Value = 558.3571000000002
$ws.Range("N113").Value = -9063.1052
$ws.Range("H131").Value = 3934.476
$ws.Range("J131").Value = 2340.375
$ws.Range("L131").Value = 7021.125
$ws.Range("N131").Value = -17101.125
$ws.Range("H137").Value = 8581
$ws.Range("J137").Value = 9876
$ws.Range("L137").Value = 29628
$ws.Range("N137").Value = -39828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 29499
$ws.Range("I18").Value = 29499
$ws.Range("K18").Value = 29499
$ws.Range("M18").Value = -29206
$ws.Range("H74").Value = 59999
$ws.Range("J74").Value = 59999
$ws.Range("L74").Value = 59999
$ws.Range("N74").Value = -61871
$ws.Range("H77").Value = 59999
$ws.Range("J77").Value = 59999
$ws.Range("L77").Value = 179997
$ws.Range("N77").Value = -189357
$ws.Range("H122").Value = 6620.4614
$ws.Range("I122").Value = 4652.25
$ws.Range("J122").Value = 9769.6
$ws.Range("K122").Value = 13956.75
$ws.Range("L122").Value = 29308.8
$ws.Range("M122").Value = -11506.75
$ws.Range("N122").Value = -34208.8
$ws.Range("H132").Value = 27780340
$ws.Range("I132").Value = 33335920
$ws.Range("K132").Value = 100007760
$ws.Range("M132").Value = -100005230

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8561.714
$ws.Range("I22").Value = 8561.714
$ws.Range("K22").Value = 8561.714
$ws.Range("M22").Value = -8266.714
$ws.Range("H27").Value = 8561.714
$ws.Range("I27").Value = 8561.714
$ws.Range("K27").Value = 8561.714
$ws.Range("M27").Value = -8454.714
$ws.Range("H46").Value = 1430.2778
$ws.Range("J46").Value = 1881.25
$ws.Range("L46").Value = 1881.25
$ws.Range("N46").Value = -2257.25
$ws.Range("H63").Value = 136000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 136000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 136000
$ws.Range("N63").Value = -137498
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 136000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 136000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 408000
$ws.Range("N66").Value = -415488
$ws.Range("M66").ClearContents()
$ws.Range("H132").Value = 28221.488
$ws.Range("I132").Value = 45429.883
$ws.Range("J132").Value = 4673.1577
$ws.Range("K132").Value = 136289.649
$ws.Range("L132").Value = 14019.4731
$ws.Range("M132").Value = -133759.649
$ws.Range("N132").Value = -19079.4731

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 57999
$ws.Range("J70").Value = 57999
$ws.Range("L70").Value = 57999
$ws.Range("N70").Value = -58629
$ws.Range("H73").Value = 57999
$ws.Range("J73").Value = 57999
$ws.Range("L73").Value = 57999
$ws.Range("N73").Value = -60183
$ws.Range("H132").Value = 771760.8
$ws.Range("I132").Value = 2769.4
$ws.Range("K132").Value = 8308.200000000001
$ws.Range("M132").Value = -5778.200000000001
